$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows 2-11 (A:id, B:nome, C:departamento, D:motivo, E:horas, F:data_serial, G:salario)
$data = @(
    @{Row=2;  A=73319; B="Dra. Bruna Ferreira";           C="TI";                     D="Outros";              E=8; F=45086; G=5032.68}
    @{Row=3;  A=57821; B="Dr. Carlos Eduardo Almeida";     C="Operações";              D="Outros";              E=6; F=45083; G=4666.58}
    @{Row=4;  A=80219; B="Mirella Pereira";                C="Atendimento ao Cliente"; D="Viagem de negócios";  E=7; F=45106; G=12055.51}
    @{Row=5;  A=18023; B="Francisco Freitas";              C="Jurídico";               D="Problemas pessoais";  E=7; F=45098; G=3277.68}
    @{Row=6;  A=1462;  B="Ian Cunha";                      C="Vendas";                 D="Outros";              E=6; F=45092; G=9471.889999999999}
    @{Row=7;  A=92505; B="Dr. Pedro Moraes";                C="Jurídico";               D="Doença";              E=6; F=45079; G=11746.11}
    @{Row=8;  A=78747; B="Vitória Dias";                   C="TI";                     D="Consulta médica";     E=5; F=45105; G=9564.16}
    @{Row=9;  A=75600; B="Luiz Miguel Caldeira";           C="Marketing";              D="Viagem de negócios";  E=2; F=45100; G=2864.58}
    @{Row=10; A=87173; B="Raquel da Rocha";                C="Operações";              D="Outros";              E=4; F=45096; G=3516.6}
    @{Row=11; A=54630; B="Pietro Rocha";                   C="Vendas";                 D="Viagem de negócios";  E=3; F=45080; G=10038.84}
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.A
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 5).Value = $item.E
    $ws.Cells.Item($r, 6).Value = $item.F
    $ws.Cells.Item($r, 7).Value = $item.G
}

$wb.Save()
